# fix(publipostage): Correct status, status name, status label ...
#
# - Drop the results_1y / results_3y / results boolean columns (old J:L)
# - Drop the old long-text "statut_name" column (old C) and fold the
#   status info into a single "statut_name" column (new B), recoding the
#   status values from the old icon/color/long-label scheme to the new
#   numeric-code scheme.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three boolean "results" columns entirely.
$ws.Columns("J:L").Delete() | Out-Null

# Remove the old verbose "statut_name" text column (old column C);
# the old "statut_label" column (old B) slides into its place and will
# be recoded below to become the new "statut_name" column.
$ws.Columns("C").Delete() | Out-Null

# Update header: old B1 was "statut_label" -> new B1 is "statut_name".
$ws.Range("B1").Value = "statut_name"

# Recode the status columns (A = statut code, B = statut_name) per row,
# based on the old A-column icon value.
$statusMap = @{
    "⚠️" = @("4", "4: pas de résultats postés ni publiés");
    "+3" = @("2", "2: résultats postés ou publiés entre 12 et 36 mois")
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $oldStatus = $cellA.Value()
    $pair = $statusMap[$oldStatus]
    if ($pair) {
        # The new "statut" codes ("4", "2") look numeric, but the source
        # data keeps them as text -- force text storage, write, then drop
        # the temporary number format again so no stray cell style lingers.
        $cellA.NumberFormat = "@"
        $cellA.Value = $pair[0]
        $cellA.ClearFormats()

        $ws.Cells.Item($r, 2).Value = $pair[1]
    }
}
